$d = $word.ActiveDocument

# The document body currently has:
#   1. An empty paragraph containing only the _GoBack bookmark
#   2-9. "Domain Knowledge" heading + body paragraphs / bullet list
# We need to remove everything except the first (bookmark) paragraph,
# leaving a single empty paragraph followed by the section properties.

$paraCount = $d.Paragraphs.Count
if ($paraCount -gt 1) {
    $startRange = $d.Paragraphs(2).Range.Start
    $endRange = $d.Paragraphs($paraCount).Range.End
    $r = $d.Range($startRange, $endRange)
    $r.Delete()
}
